$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.271.71"
$ws.Range("E2").Value = "  +0.29%  "

# Row 3
$ws.Range("D3").Value = "1.590.37"
$ws.Range("E3").Value = "  +0.52%  "

# Row 4
$ws.Range("E4").Value = "  -0.22%  "

# Row 5
$ws.Range("D5").Value = "'213.19"
$ws.Range("E5").Value = "  +1.66%  "

# Row 6
$ws.Range("D6").Value = "'0.500"
$ws.Range("E6").Value = "  +0.68%  "

# Row 7
$ws.Range("E7").Value = "  -0.23%  "

# Row 8
$ws.Range("E8").Value = "  +0.20%  "

# Row 9
$ws.Range("E9").Value = "  -0.19%  "

# Row 10
$ws.Range("E10").Value = "  -0.85%  "

# Row 11
$ws.Range("D11").Value = "'0.0848"

# Row 12
$ws.Range("D12").Value = "1.813.34"
$ws.Range("E12").Value = "  +0.46%  "

# Row 13
$ws.Range("D13").Value = "1.612.72"
$ws.Range("E13").Value = "  +1.68%  "

# Row 14
$ws.Range("E14").Value = "  -0.15%  "

# Row 15
$ws.Range("E15").Value = "  +1.41%  "

# Row 16
$ws.Range("D16").Value = "'64.49"
$ws.Range("E16").Value = "  +0.01%  "

# Row 17
$ws.Range("D17").Value = "26.280.41"
$ws.Range("E17").Value = "  +0.31%  "

# Row 18
$ws.Range("E18").Value = "  -0.82%  "

# Row 19
$ws.Range("D19").Value = "'7.47"
$ws.Range("E19").Value = "  +2.71%  "

# Row 20
$ws.Range("D20").Value = "'213.31"
$ws.Range("E20").Value = "  +2.98%  "

# Row 21
$ws.Range("E21").Value = "  -0.18%  "

# Row 22
$ws.Range("E22").Value = "  +0.80%  "

# Row 23
$ws.Range("B23").Value = "Avalanche"
$ws.Range("C23").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D23").Value = "'8.96"
$ws.Range("E23").Value = "  +0.82%  "

# Row 24
$ws.Range("B24").Value = "Toncoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D24").Value = "'2.15"
$ws.Range("E24").Value = "  -2.13%  "

# Row 25
$ws.Range("D25").Value = "'145.00"
$ws.Range("E25").Value = "  +0.27%  "

# Row 27
$ws.Range("E27").Value = "  +0.79%  "

# Row 28
$ws.Range("E28").Value = "  -0.34%  "

# Row 29
$ws.Range("E29").Value = "  -0.17%  "

# Row 30
$ws.Range("D30").Value = "'0.0499"
$ws.Range("E30").Value = "  -0.88%  "

# Row 31
$ws.Range("E31").Value = "  +1.03%  "

# Row 32
$ws.Range("E32").Value = "  -0.03%  "

# Row 33
$ws.Range("E33").Value = "  +0.37%  "

# Row 34
$ws.Range("D34").Value = "1.338.77"
$ws.Range("E34").Value = "  +4.89%  "

# Row 35
$ws.Range("E35").Value = "  -0.90%  "

# Row 36
$ws.Range("E36").Value = "  -0.34%  "

# Row 37
$ws.Range("D37").Value = "'0.593"
$ws.Range("E37").Value = "  -3.02%  "

# Row 38
$ws.Range("E38").Value = "  -0.20%  "

# Row 39
$ws.Range("E39").Value = "  +0.05%  "

# Row 40
$ws.Range("D40").Value = "'5.78"
$ws.Range("E40").Value = "  +4.16%  "

# Row 41
$ws.Range("E41").Value = "  -0.19%  "

# Row 42
$ws.Range("E42").Value = "  -0.44%  "

# Row 43
$ws.Range("E43").Value = "  +0.43%  "

# Row 44
$ws.Range("E44").Value = "  -0.44%  "

# Row 45
$ws.Range("D45").Value = "1.724.89"
$ws.Range("E45").Value = "  +0.29%  "

# Row 46
$ws.Range("D46").Value = "'61.81"
$ws.Range("E46").Value = "  -0.76%  "

# Row 47
$ws.Range("D47").Value = "'87.23"
$ws.Range("E47").Value = "  -2.08%  "

# Row 48
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "0.0₆0104"
$ws.Range("E48").Value = "  -1.23%  "

# Row 49
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").Value = "'1.50"
$ws.Range("E49").Value = "  -3.73%  "

# Row 50
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "'0.0504"
$ws.Range("E50").Value = "  -0.50%  "

# Row 51
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "'0.0979"
$ws.Range("E51").Value = "  -2.39%  "
